$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "[Rogério-Retífica-2NB, -, -, -]"
$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"
$ws.Range("B4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("B6").Value = "[Rogério-Retífica-2NB, -, -, -]"
$ws.Range("F6").Value = "-"
$ws.Range("B7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "-"
$ws.Range("C8").Value = "[Rogério-Retífica-2NB, -, -, -]"
$ws.Range("D8").Value = "[Rogério-Retífica-2NB, -, -, -]"
$ws.Range("F8").Value = "-"
$ws.Range("B10").Value = "-"
$ws.Range("D10").Value = "[Rogério-Retífica-2NB, -, -, -]"
$ws.Range("E10").Value = "-"
$ws.Range("B11").Value = "-"
$ws.Range("F11").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("F14").Value = "[Rogério-Retífica-2NB, -, -, -]"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "-"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("B18").Value = "Suzanny-Des. Maq. Cad-"
$ws.Range("C18").Value = "[Joel L.-Fundição-2NB, Elcio Dec.-C.pneumática-2NB, Ismail-Metrologia 2-2NB, Cláudio-Soldagem-2NB]"
$ws.Range("D18").Value = "[Ismail-Metrologia 2-2NB, Paulo Rob.-CAD/CAM-2NB, Leandro-M.S.R.A.C.-2NB, Victor-Usin. CNC-2NB]"
$ws.Range("E18").Value = "[Leandro-M.S.R.A.C.-2NB, Victor-Usin. CNC-2NB, Leandro-M. Maq. E. I.-2NB, Paulo Rob.-CAD/CAM-2NB]"
$ws.Range("F18").Value = "[Rogério-Retífica-2NB, Guilherme-C. L. P.-2NB, Aderci-Fresagem-2NB, Guilherme-Coman. Hidraulicos-2NB]"
$ws.Range("B19").Value = "Suzanny-Des. Maq. Cad-"
$ws.Range("C19").Value = "[Rogério-Retífica-2NB, Elcio Dec.-C.pneumática-2NB, Ismail-Metrologia 2-2NB, Cláudio-Soldagem-2NB]"
$ws.Range("D19").Value = "[Leandro-M. Maq. E. I.-2NB, Paulo Rob.-CAD/CAM-2NB, Leandro-M.S.R.A.C.-2NB, Victor-Usin. CNC-2NB]"
$ws.Range("E19").Value = "Claudinei-Elemaq.-"
$ws.Range("F19").Value = "[Rogério-Retífica-2NB, Guilherme-C. L. P.-2NB, Aderci-Fresagem-2NB, Guilherme-Coman. Hidraulicos-2NB]"
$ws.Range("B20").Value = "Suzanny-Des. Maq. Cad-"
$ws.Range("C20").Value = "[Joel L.-Fundição-2NB, Elcio Dec.-C.pneumática-2NB, Ismail-Metrologia 2-2NB, Cláudio-Soldagem-2NB]"
$ws.Range("D20").Value = "[Leandro-M. Maq. E. I.-2NB, Paulo Rob.-CAD/CAM-2NB, Leandro-M.S.R.A.C.-2NB, Victor-Usin. CNC-2NB]"
$ws.Range("E20").Value = "Claudinei-Elemaq.-"
$ws.Range("F20").Value = "[Joel L.-Fundição-2NB, Guilherme-C. L. P.-2NB, Aderci-Fresagem-2NB, Guilherme-Coman. Hidraulicos-2NB]"
$ws.Range("B21").Value = "Euclides-Gestão integrada-"
$ws.Range("C21").Value = "[Joel L.-Fundição-2NB, Elcio Dec.-C.pneumática-2NB, Leandro-M. Maq. E. I.-2NB, Cláudio-Soldagem-2NB]"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "Euclides-Gestão integrada-"
$ws.Range("F21").Value = "[Rogério-Retífica-2NB, Guilherme-C. L. P.-2NB, Aderci-Fresagem-2NB, Guilherme-Coman. Hidraulicos-2NB]"
